$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 192-196: the match between "Junior" and "Huila" (row 194) stays put,
#    but the two match-pairs around it got re-ordered by the scraper re-run.
#    Swap the match-data columns (F..V, minus the timestamp columns that are
#    identical between the two rows anyway) between row 192 <-> row 196 and
#    row 193 <-> row 195. Column A (Indice) and the opening-odds timestamp
#    columns K/O/S are untouched because they were already identical.
# ---------------------------------------------------------------------------
$swapCols = @("F","G","H","J","L","M","N","P","Q","R","T","U","V")

function Swap-RowData($r1, $r2) {
    foreach ($col in $swapCols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

Swap-RowData 192 196
Swap-RowData 193 195

# ---------------------------------------------------------------------------
# 2) Append four newly scraped matches as rows 206-209.
# ---------------------------------------------------------------------------
function Add-MatchRow($m) {
    $Row = $m.Row

    $ws.Range("A$Row").Value = $m.Indice
    $ws.Range("B$Row").Value = "colombia"
    $ws.Range("C$Row").Value = "primera-a"
    $ws.Range("D$Row").Value = "'2023"
    $ws.Range("E$Row").Value = $m.DataPartida
    $ws.Range("F$Row").Value = $m.Home
    $ws.Range("G$Row").Value = $m.HomeGols
    $ws.Range("H$Row").Value = $m.Away
    $ws.Range("I$Row").Value = $m.AwayGols
    $ws.Range("J$Row").Value = $m.HomeOpenOdds
    $ws.Range("K$Row").Value = $m.HomeOpenDt
    $ws.Range("L$Row").Value = $m.HomeCloseOdds
    $ws.Range("M$Row").Value = $m.HomeCloseDt
    $ws.Range("N$Row").Value = $m.DrawOpenOdds
    $ws.Range("O$Row").Value = $m.DrawOpenDt
    $ws.Range("P$Row").Value = $m.DrawCloseOdds
    $ws.Range("Q$Row").Value = $m.DrawCloseDt
    $ws.Range("R$Row").Value = $m.AwayOpenOdds
    $ws.Range("S$Row").Value = $m.AwayOpenDt
    $ws.Range("T$Row").Value = $m.AwayCloseOdds
    $ws.Range("U$Row").Value = $m.AwayCloseDt
    $ws.Range("V$Row").Value = $m.Url

    # Match the style used by the rest of the table: column A is bold with a
    # border (style carried over from the header row), column E is the
    # date/time formatted serial number.
    $ws.Range("A205").Copy() | Out-Null
    $ws.Range("A$Row").PasteSpecial(-4122) | Out-Null
    $ws.Range("E205").Copy() | Out-Null
    $ws.Range("E$Row").PasteSpecial(-4122) | Out-Null
}

Add-MatchRow @{
    Row = 206; Indice = 205; DataPartida = 45248.95833333334
    Home = "Aguilas"; HomeGols = 0; Away = "Deportes Tolima"; AwayGols = 4
    HomeOpenOdds = 2.06; HomeOpenDt = "14/11/2023 02:12"
    HomeCloseOdds = 1.96; HomeCloseDt = "18/11/2023 22:59"
    DrawOpenOdds = 3.26; DrawOpenDt = "14/11/2023 02:12"
    DrawCloseOdds = 3.23; DrawCloseDt = "18/11/2023 22:59"
    AwayOpenOdds = 3.97; AwayOpenDt = "14/11/2023 02:12"
    AwayCloseOdds = 4.63; AwayCloseDt = "18/11/2023 22:59"
    Url = "https://www.betexplorer.com/football/colombia/primera-a/aguilas-doradas-deportes-tolima/YoMfV0XC/"
}

Add-MatchRow @{
    Row = 207; Indice = 206; DataPartida = 45249.0625
    Home = "Junior"; HomeGols = 3; Away = "Dep. Cali"; AwayGols = 0
    HomeOpenOdds = 1.75; HomeOpenDt = "14/11/2023 02:12"
    HomeCloseOdds = 1.74; HomeCloseDt = "19/11/2023 01:24"
    DrawOpenOdds = 3.53; DrawOpenDt = "14/11/2023 02:12"
    DrawCloseOdds = 3.59; DrawCloseDt = "19/11/2023 01:24"
    AwayOpenOdds = 5.18; AwayOpenDt = "14/11/2023 02:12"
    AwayCloseOdds = 5.45; AwayCloseDt = "19/11/2023 01:24"
    Url = "https://www.betexplorer.com/football/colombia/primera-a/junior-dep-cali/WWvbUKnJ/"
}

Add-MatchRow @{
    Row = 208; Indice = 207; DataPartida = 45249.91666666666
    Home = "Millonarios"; HomeGols = 2; Away = "America De Cali"; AwayGols = 1
    HomeOpenOdds = 2.03; HomeOpenDt = "13/11/2023 01:42"
    HomeCloseOdds = 1.99; HomeCloseDt = "19/11/2023 21:54"
    DrawOpenOdds = 3.3; DrawOpenDt = "13/11/2023 01:42"
    DrawCloseOdds = 3.25; DrawCloseDt = "19/11/2023 21:54"
    AwayOpenOdds = 4.03; AwayOpenDt = "13/11/2023 01:42"
    AwayCloseOdds = 4.45; AwayCloseDt = "19/11/2023 21:54"
    Url = "https://www.betexplorer.com/football/colombia/primera-a/millonarios-america-de-cali/WAs6xIvt/"
}

Add-MatchRow @{
    Row = 209; Indice = 208; DataPartida = 45250.02083333334
    Home = "Ind. Medellin"; HomeGols = 2; Away = "Atl. Nacional"; AwayGols = 1
    HomeOpenOdds = 1.83; HomeOpenDt = "13/11/2023 01:42"
    HomeCloseOdds = 2.04; HomeCloseDt = "20/11/2023 00:22"
    DrawOpenOdds = 3.59; DrawOpenDt = "13/11/2023 01:42"
    DrawCloseOdds = 3.39; DrawCloseDt = "20/11/2023 00:22"
    AwayOpenOdds = 4.58; AwayOpenDt = "13/11/2023 01:42"
    AwayCloseOdds = 4.02; AwayCloseDt = "20/11/2023 00:22"
    Url = "https://www.betexplorer.com/football/colombia/primera-a/ind-medellin-atl-nacional/A1tAyxfn/"
}

# ---------------------------------------------------------------------------
# 3) Keep the sheet's declared dimension/used-range in sync (A1:V205 -> A1:V209).
# ---------------------------------------------------------------------------
$ws.Range("A1:V209").Select() | Out-Null
